$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 4 name first
$ws.Cells.Item(4, 1).Value = "almuerzo con DELL"

# event_type column (B) for the new rows
$ws.Cells.Item(4, 2).Value = "O"
$ws.Cells.Item(5, 2).Value = "F"
$ws.Cells.Item(6, 2).Value = "O"
$ws.Cells.Item(7, 2).Value = "F"
$ws.Cells.Item(8, 2).Value = "F"

# sponsor_name column (E) filled top to bottom
$ws.Cells.Item(4, 5).Value = "DELL"
$ws.Cells.Item(5, 5).Value = "ESTEBAN GUTIERREZ"
$ws.Cells.Item(6, 5).Value = "PEPSI"
$ws.Cells.Item(7, 5).Value = "ANDRES GOMEZ"
$ws.Cells.Item(8, 5).Value = "MICROSOFT"

# participation column (F) filled top to bottom
$ws.Cells.Item(4, 6).Value = "patrocinador del almuerzo"
$ws.Cells.Item(5, 6).Value = "pago para estudientes"
$ws.Cells.Item(6, 6).Value = "revicion de notas estudientes becados"
$ws.Cells.Item(7, 6).Value = "desayuno con Andres Gomez"
$ws.Cells.Item(8, 6).Value = "reparto de equipos"

# Remaining name column (A) entries, filled out of order
$ws.Cells.Item(7, 1).Value = "desayuno con Andres Gomez"
$ws.Cells.Item(8, 1).Value = "entrega de material Microsoft"
$ws.Cells.Item(6, 1).Value = "reunion de control estudientes"
$ws.Cells.Item(5, 1).Value = "recaudo becas"

$ws.Range("A6").Select()
